# Table A2 - Descriptives WM Performance
# Reformat the CL column values ("0.00"/"1.00" -> "0"/"1") and swap the
# "high"/"low" rows' data within each CL group.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Exclude the trailing cell-end mark so we don't clobber the cell
    # structure - only replace the visible text.
    $r.End = $r.End - 1
    $r.Text = $text
}

# Row 3: CL 0.00 / high  ->  CL 0 / low
Set-CellText $t 3 1 "0"
Set-CellText $t 3 2 "low"
Set-CellText $t 3 3 "0.63"
Set-CellText $t 3 4 "0.14"
Set-CellText $t 3 5 "3.56"
Set-CellText $t 3 6 "2.45"

# Row 4: CL 0.00 / low  ->  CL 0 / high
Set-CellText $t 4 1 "0"
Set-CellText $t 4 2 "high"
Set-CellText $t 4 3 "0.73"
Set-CellText $t 4 4 "0.16"
Set-CellText $t 4 5 "9.44"
Set-CellText $t 4 6 "4.31"

# Row 5: CL 1.00 / high  ->  CL 1 / low
Set-CellText $t 5 1 "1"
Set-CellText $t 5 2 "low"
Set-CellText $t 5 3 "0.61"
Set-CellText $t 5 4 "0.11"
Set-CellText $t 5 5 "4.28"
Set-CellText $t 5 6 "3.79"

# Row 6: CL 1.00 / low  ->  CL 1 / high
Set-CellText $t 6 1 "1"
Set-CellText $t 6 2 "high"
Set-CellText $t 6 3 "0.69"
Set-CellText $t 6 4 "0.16"
Set-CellText $t 6 5 "10.21"
Set-CellText $t 6 6 "5.12"
